$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a "price" cell (column D) as literal text, preventing Excel's
# automatic number/date coercion from mangling values like "141.30",
# "57.461.82" or "0.0₃0927". A leading apostrophe forces text entry, exactly
# like a user typing '141.30 into the cell.
function Set-PriceText {
    param($addr, $text)
    $ws.Range($addr).Value = "'" + $text
}

# Row 2 - Bitcoin
Set-PriceText "D2" "57.461.82"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3 - Ethereum
Set-PriceText "D3" "3.084.61"
$ws.Range("E3").Value = "  +0.65%  "

# Row 5 - BNB
Set-PriceText "D5" "515.51"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6 - Solana
Set-PriceText "D6" "141.30"
$ws.Range("E6").Value = "  +0.48%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.12%  "

# Row 8 - XRP
Set-PriceText "D8" "0.434"
$ws.Range("E8").Value = "  -0.37%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  +0.36%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.19%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -1.30%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-PriceText "D12" "3.616.46"
$ws.Range("E12").Value = "  +0.88%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +2.61%  "

# Row 14 - Avalanche
Set-PriceText "D14" "25.71"
$ws.Range("E14").Value = "  -4.42%  "

# Row 15 - ShibaInu
Set-PriceText "D15" "0.0000164"
$ws.Range("E15").Value = "  -2.70%  "

# Row 16 - WrappedBTC
Set-PriceText "D16" "57.579.21"
$ws.Range("E16").Value = "  -0.27%  "

# Row 17 - WrappedEther
Set-PriceText "D17" "3.085.74"
$ws.Range("E17").Value = "  +0.73%  "

# Row 18 - Polkadot
Set-PriceText "D18" "6.14"
$ws.Range("E18").Value = "  -1.15%  "

# Row 19 - Chainlink
Set-PriceText "D19" "13.12"
$ws.Range("E19").Value = "  -2.90%  "

# Row 20 - Uniswap
Set-PriceText "D20" "8.15"
$ws.Range("E20").Value = "  -0.58%  "

# Row 21 - BitcoinCash
Set-PriceText "D21" "335.38"
$ws.Range("E21").Value = "  +1.01%  "

# Row 22 - Dai
Set-PriceText "D22" "1.00"
$ws.Range("E22").Value = "  +0.14%  "

# Row 23 - Polygon
Set-PriceText "D23" "0.502"
$ws.Range("E23").Value = "  -1.34%  "

# Row 24 - Litecoin
Set-PriceText "D24" "65.84"
$ws.Range("E24").Value = "  +0.83%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  +3.75%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.21%  "

# Row 27 - PEPE
Set-PriceText "D27" "0.0₃0927"
$ws.Range("E27").Value = "  +1.65%  "

# Row 28 - RenderToken
Set-PriceText "D28" "6.44"
$ws.Range("E28").Value = "  -4.97%  "

# Row 29 - InternetComputer(DFINITY)
Set-PriceText "D29" "7.13"
$ws.Range("E29").Value = "  -1.93%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.19%  "

# Row 31 - EthereumClassic
Set-PriceText "D31" "20.86"
$ws.Range("E31").Value = "  -0.31%  "

# Row 32 - Fetch.AI
Set-PriceText "D32" "1.17"
$ws.Range("E32").Value = "  -3.96%  "

# Row 33 - Monero
Set-PriceText "D33" "153.99"
$ws.Range("E33").Value = "  +0.03%  "

# Row 34 - was EnergySwap, now NEARProtocol
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-PriceText "D34" "4.55"
$ws.Range("E34").Value = "  -3.07%  "

# Row 35 - was NEARProtocol, now EnergySwap
$ws.Range("B35").Value = "EnergySwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-PriceText "D35" "27.65"
$ws.Range("E35").Value = "  +10.08%  "

# Row 36 - Aptos
Set-PriceText "D36" "5.93"
$ws.Range("E36").Value = "  -0.08%  "

# Row 37 - ImmutableX
Set-PriceText "D37" "1.25"
$ws.Range("E37").Value = "  -1.90%  "

# Row 38 - Hedera
Set-PriceText "D38" "0.0678"
$ws.Range("E38").Value = "  -0.61%  "

# Row 39 - RenzoRestakedETH
Set-PriceText "D39" "3.125.28"
$ws.Range("E39").Value = "  +0.93%  "

# Row 40 - OKB
Set-PriceText "D40" "36.81"
$ws.Range("E40").Value = "  -1.10%  "

# Row 41 - Mantle
$ws.Range("E41").Value = "  +0.37%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  -1.49%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.14%  "

# Row 44 - Maker
Set-PriceText "D44" "2.292.17"
$ws.Range("E44").Value = "  +3.72%  "

# Row 45 - VeChain
Set-PriceText "D45" "0.0254"
$ws.Range("E45").Value = "  +4.23%  "

# Row 46 - Stacks
$ws.Range("E46").Value = "  -0.70%  "

# Row 47 - ONDO
Set-PriceText "D47" "0.946"
$ws.Range("E47").Value = "  -0.72%  "

# Row 48 - InjectiveProtocol
Set-PriceText "D48" "20.14"
$ws.Range("E48").Value = "  -0.69%  "

# Row 49 - Cosmos
Set-PriceText "D49" "5.90"
$ws.Range("E49").Value = "  -3.51%  "

# Row 51 - Bittensor
Set-PriceText "D51" "248.80"
$ws.Range("E51").Value = "  +7.19%  "
